$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value, derived from the day-over-day
# cryptos price/volume refresh produced by the scheduled GitHub Action.
$updates = [ordered]@{
    "D2" = "65.738.45"
    "E2" = "  -0.02%  "
    "D3" = "2.662.25"
    "E3" = "  -0.58%  "
    "E4" = "  +0.06%  "
    "D5" = "598.89"
    "E5" = "  -0.32%  "
    "D6" = "159.70"
    "E6" = "  +1.89%  "
    "D7" = "0.641"
    "E7" = "  +3.82%  "
    "D8" = "1.00"
    "E8" = "  +0.06%  "
    "D9" = "0.127"
    "E9" = "  -1.95%  "
    "E10" = "  -0.39%  "
    "E11" = "  -0.18%  "
    "E12" = "  +1.22%  "
    "D13" = "29.23"
    "E13" = "  -0.20%  "
    "D14" = "0.0000195"
    "E14" = "  -1.35%  "
    "D15" = "3.142.34"
    "E15" = "  -0.48%  "
    "D16" = "65.692.98"
    "E16" = "  +0.13%  "
    "D17" = "2.677.43"
    "E17" = "  +0.07%  "
    "D18" = "12.56"
    "E18" = "  -2.74%  "
    "E19" = "  -0.07%  "
    "D20" = "354.14"
    "E20" = "  +0.39%  "
    "D21" = "7.47"
    "E21" = "  -1.41%  "
    "D22" = "1.00"
    "E22" = "  -0.04%  "
    "D23" = "69.82"
    "E23" = "  -0.07%  "
    "D24" = "1.78"
    "E24" = "  +8.10%  "
    "E25" = "  +0.44%  "
    "D26" = "9.75"
    "E26" = "  +1.12%  "
    "E27" = "  +1.94%  "
    "D28" = "563.56"
    "E28" = "  +5.54%  "
    "D29" = "8.14"
    "E29" = "  +0.99%  "
    "E30" = "  -2.20%  "
    "E31" = "  +0.07%  "
    "E32" = "  +0.26%  "
    "D33" = "1.82"
    "E33" = "  +2.96%  "
    "E34" = "  +3.36%  "
    "E35" = "  -0.43%  "
    "E36" = "  -0.16%  "
    "D37" = "20.62"
    "E37" = "  -0.01%  "
    "E38" = "  +1.86%  "
    "E39" = "  +0.01%  "
    "D40" = "154.44"
    "E40" = "  -3.17%  "
    "E41" = "  +7.41%  "
    "D42" = "161.55"
    "E42" = "  -1.36%  "
    "E43" = "  -0.84%  "
    "D44" = "0.0618"
    "E44" = "  +1.22%  "
    "D45" = "23.46"
    "E45" = "  +2.81%  "
    "D46" = "0.646"
    "E46" = "  +0.90%  "
    "E47" = "  +0.23%  "
    "E48" = "  +1.86%  "
    "D49" = "19.89"
    "E49" = "  -1.89%  "
    "D50" = "0.0₆0245"
    "E50" = "  -6.97%  "
    "D51" = "0.817"
    "E51" = "  +0.51%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "1.00", "159.70")
    # are not silently coerced into numbers, matching the source data.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
